# 完成了"缩放不流畅"和"加载慢"问题的记录更新
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 13: mark as resolved ----
# C13: add description of the fix (don't need the thumbnail image, show live instead)
$c = $ws.Range("C13")
$c.Value = "不需要这个图片，要进行实时显示"
$c.WrapText = $true
$c.VerticalAlignment = -4108

# E13: mark solved status as OK
$e = $ws.Range("E13")
$e.Value = "OK"
$e.WrapText = $true
$e.HorizontalAlignment = -4108
$e.VerticalAlignment = -4108

# ---- Row 27: finish the "load slow" note, fill in solution / status / date ----
$ws.Range("C27").Value = "估计是内存的占用问题，图片加载效率还是太低了。"

$d27 = $ws.Range("D27")
$d27.Value = "在每次点击item的时候，一直进行notifyDataSetChanged，这个严重影响了刷新的速度。去掉之后，在图片界面显示就比较快了。"
$d27.WrapText = $true
$d27.VerticalAlignment = -4108

$e27 = $ws.Range("E27")
$e27.Value = "ok"
$e27.WrapText = $true
$e27.VerticalAlignment = -4108

$f27 = $ws.Range("F27")
$f27.Value = "2016.5.10"
$f27.WrapText = $true
$f27.VerticalAlignment = -4108

$ws.Rows.Item(27).RowHeight = 45

# ---- New row 28: gridview load-folder picture issue ----
$a28 = $ws.Range("A28"); $a28.Value = 28; $a28.WrapText = $true; $a28.VerticalAlignment = -4108
$b28 = $ws.Range("B28"); $b28.Value = "加载指定文件夹的图片"; $b28.WrapText = $true; $b28.VerticalAlignment = -4108
$d28 = $ws.Range("D28"); $d28.Value = "对查询出出来的图片进行判断"; $d28.WrapText = $true; $d28.VerticalAlignment = -4108
$e28 = $ws.Range("E28"); $e28.Value = "ok"; $e28.WrapText = $true; $e28.VerticalAlignment = -4108
$f28 = $ws.Range("F28"); $f28.Value = "2016.5.10"; $f28.WrapText = $true; $f28.VerticalAlignment = -4108

# ---- New row 29: gridview setSelection issue ----
$a29 = $ws.Range("A29"); $a29.Value = 29; $a29.WrapText = $true; $a29.VerticalAlignment = -4108
$b29 = $ws.Range("B29"); $b29.Value = "gridview的setselection无效"; $b29.WrapText = $true; $b29.VerticalAlignment = -4108
$d29 = $ws.Range("D29"); $d29.Value = "莫名其妙就好了，也许真不应该加上notifyDataSetChanged。"; $d29.WrapText = $true; $d29.VerticalAlignment = -4108
$e29 = $ws.Range("E29"); $e29.Value = "ok"; $e29.WrapText = $true; $e29.VerticalAlignment = -4108
$f29 = $ws.Range("F29"); $f29.Value = "2016.5.10"; $f29.WrapText = $true; $f29.VerticalAlignment = -4108

# ---- New row 30: zoom not smooth issue ----
$a30 = $ws.Range("A30"); $a30.Value = 30; $a30.WrapText = $true; $a30.VerticalAlignment = -4108
$b30 = $ws.Range("B30"); $b30.Value = "大图的缩放问题，图片显示的不正常"; $b30.WrapText = $true; $b30.VerticalAlignment = -4108
$d30 = $ws.Range("D30"); $d30.Value = "matrix图形变换，实现图片的缩放，效果比较好了"; $d30.WrapText = $true; $d30.VerticalAlignment = -4108
$e30 = $ws.Range("E30"); $e30.Value = "ok"; $e30.WrapText = $true; $e30.VerticalAlignment = -4108
$f30 = $ws.Range("F30"); $f30.Value = "2016.5.10"; $f30.WrapText = $true; $f30.VerticalAlignment = -4108

# ---- New row 31: pause/play icon bug on back-press ----
$a31 = $ws.Range("A31"); $a31.Value = 31; $a31.WrapText = $true; $a31.VerticalAlignment = -4108
$b31 = $ws.Range("B31"); $b31.Value = "音乐暂停状态下，按下返回键，播放按钮上的暂停图片变成了播放图片"; $b31.WrapText = $true; $b31.VerticalAlignment = -4108
$ws.Range("C31").Value = "onbackpress函数中修改"
$e31 = $ws.Range("E31"); $e31.Value = "ok"; $e31.WrapText = $true; $e31.VerticalAlignment = -4108
$f31 = $ws.Range("F31"); $f31.Value = "2016.5.10"; $f31.WrapText = $true; $f31.VerticalAlignment = -4108

# ---- column F gets an explicit width ----
$ws.Columns.Item(6).ColumnWidth = 12.4

# ---- move the frozen-pane view / selection down to the new last rows ----
$ws.Range("D28").Select()
